$d = $word.ActiveDocument

# --- Paragraph 1: split "In Module " into "In " + a long narrative that
# wraps around the existing _GoBack bookmark. --------------------------

# 1) Shrink the original run's text down to "In " (the _GoBack bookmark
#    still sits immediately after it, at the end of the paragraph). Using
#    a range (rather than $d.Content directly) lets us read back exactly
#    where the replaced text ended, instead of hard-coding its length.
$findRange = $d.Content
$findRange.Find.Execute("In Module ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "In ", 2) | Out-Null
$bookmarkPos = $findRange.End

# 2) The text that must land *after* the bookmark can't be inserted
#    directly at the bookmark's position (new text at a trailing
#    zero-width bookmark always lands before it). Work around this by
#    appending a new paragraph at the very end of the document, putting
#    that text there, then deleting the paragraph mark that separates it
#    from paragraph 1 - this merges it back in as a new run positioned
#    after the bookmark, inside the same paragraph.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter() | Out-Null

$afterBookmarkPara = $d.Paragraphs(2)
$afterBookmarkPara.Range.InsertAfter( `
    " how to do this prior to talking with Ryan. I was able to also become a member of the repository that was created for our project. ") | Out-Null

$para1 = $d.Paragraphs(1)
$joinMark = $d.Range($para1.Range.End - 1, $para1.Range.End)
$joinMark.Delete() | Out-Null

# 3) Insert the long lead-in text before the bookmark. A range collapsed
#    to the bookmark's position always receives InsertBefore text
#    immediately in front of the bookmark, so inserting chunk-by-chunk
#    (collapsing to the end after each one) builds up the paragraph in
#    the correct left-to-right order.
$before = $d.Range($bookmarkPos, $bookmarkPos)
$before.InsertBefore("module 2, I learned how to create a use case for our project. I worked with my team to determine the type of use case that I was going to create and submit. I ") | Out-Null
$before.Collapse(0)
$before.InsertBefore("chose to create a use case that described how a user account would be deleted from our website/application. Upon developing this, ") | Out-Null
$before.Collapse(0)
$before.InsertBefore("I learned what exactly the use case was. I also created a diagram that showed the flow of our application and how it was going to work. I worked with Ryan to ") | Out-Null
$before.Collapse(0)
$before.InsertBefore("commit my use case document to my own branch within our project. I then created a pull request to be reviewed by my team in order to merge my branch with the master branch. ") | Out-Null
$before.Collapse(0)
$before.InsertBefore("I thought this was interesting because I was not sure") | Out-Null

# --- New trailing paragraphs -------------------------------------------

# 4) Empty paragraph right after paragraph 1.
$emptyParaRange = $d.Content
$emptyParaRange.Collapse(0)
$emptyParaRange.Text = [char]13

# 5) Final paragraph with the Slack sentence.
$slackParaRange = $d.Content
$slackParaRange.Collapse(0)
$slackParaRange.Text = [char]13

$slackTextRange = $d.Content
$slackTextRange.Collapse(0)
$slackTextRange.InsertBefore("I was also able to communicate on Slack with my team regarding information for our project. ") | Out-Null
